# Applies the "Updated cryptos list on Tue Jul  4 17:44:14 UTC 2023 with GitHub Actions"
# commit: refreshed Price (D) / Volume(1h) (E) figures for each coin row, and
# restored rows 18/19 (Avalanche/ShibaInu) and 46/47 (Aptos/PaxDollar) to their
# correct rank order with refreshed data.
#
# Cells are forced to Text ("@") format before the write so numeric-looking
# strings (e.g. "1.002", "0.000007744") are stored as literal text, matching
# the source feed's inlineStr cells, then the style is reset to Normal so we
# do not leave a stray cell format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "30.959.07"
Set-TextValue "E2" "  -0.34%  "

Set-TextValue "D3" "1.954.56"
Set-TextValue "E3" "  -0.64%  "

Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.12%  "

Set-TextValue "D5" "243.11"
Set-TextValue "E5" "  -2.12%  "

Set-TextValue "E6" "  +0.00%  "

Set-TextValue "D7" "0.4860"
Set-TextValue "E7" "  -0.37%  "

Set-TextValue "D8" "0.2932"
Set-TextValue "E8" "  -0.93%  "

Set-TextValue "D9" "0.07011"
Set-TextValue "E9" "  +2.63%  "

Set-TextValue "D10" "19.45"
Set-TextValue "E10" "  +1.19%  "

Set-TextValue "D11" "107.45"
Set-TextValue "E11" "  -0.04%  "

Set-TextValue "D12" "1.949.81"
Set-TextValue "E12" "  -0.86%  "

Set-TextValue "D13" "0.07751"
Set-TextValue "E13" "  -0.41%  "

Set-TextValue "D14" "5.351"
Set-TextValue "E14" "  -1.84%  "

Set-TextValue "D15" "0.6980"
Set-TextValue "E15" "  -0.69%  "

Set-TextValue "D16" "277.55"
Set-TextValue "E16" "  -3.44%  "

Set-TextValue "D17" "30.974.26"
Set-TextValue "E17" "  -0.32%  "

Set-TextValue "B18" "ShibaInu"
Set-TextValue "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D18" "0.000007744"
Set-TextValue "E18" "  -0.08%  "

Set-TextValue "B19" "Avalanche"
Set-TextValue "C19" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D19" "13.21"
Set-TextValue "E19" "  -0.23%  "

Set-TextValue "D20" "2.209.88"
Set-TextValue "E20" "  -0.62%  "

Set-TextValue "E21" "  -0.07%  "

Set-TextValue "D22" "5.469"
Set-TextValue "E22" "  -2.90%  "

Set-TextValue "D23" "1.003"
Set-TextValue "E23" "  -0.15%  "

Set-TextValue "D24" "6.484"
Set-TextValue "E24" "  -2.07%  "

Set-TextValue "D25" "9.738"
Set-TextValue "E25" "  -2.87%  "

Set-TextValue "D26" "168.52"
Set-TextValue "E26" "  -1.26%  "

Set-TextValue "D27" "19.64"
Set-TextValue "E27" "  -2.14%  "

Set-TextValue "D28" "2.163"
Set-TextValue "E28" "  -1.72%  "

Set-TextValue "D29" "0.1045"
Set-TextValue "E29" "  -2.40%  "

Set-TextValue "D30" "1.403"
Set-TextValue "E30" "  -3.14%  "

Set-TextValue "D31" "1.562"
Set-TextValue "E31" "  -2.55%  "

Set-TextValue "D32" "4.583"
Set-TextValue "E32" "  -5.10%  "

Set-TextValue "D33" "4.390"
Set-TextValue "E33" "  -2.76%  "

Set-TextValue "D34" "0.04872"
Set-TextValue "E34" "  -4.54%  "

Set-TextValue "D35" "0.7504"
Set-TextValue "E35" "  -3.10%  "

Set-TextValue "D36" "1.162"
Set-TextValue "E36" "  -0.92%  "

Set-TextValue "D37" "2.732"
Set-TextValue "E37" "  +0.04%  "

Set-TextValue "D38" "0.01992"
Set-TextValue "E38" "  -2.61%  "

Set-TextValue "D39" "2.678"
Set-TextValue "E39" "  -1.83%  "

Set-TextValue "D40" "6.518"
Set-TextValue "E40" "  +0.22%  "

Set-TextValue "D41" "78.06"
Set-TextValue "E41" "  +6.90%  "

Set-TextValue "D42" "2.104"
Set-TextValue "E42" "  -1.49%  "

Set-TextValue "D43" "0.8951"
Set-TextValue "E43" "  +0.43%  "

Set-TextValue "D44" "109.17"
Set-TextValue "E44" "  -0.69%  "

Set-TextValue "D45" "0.4427"
Set-TextValue "E45" "  -1.66%  "

Set-TextValue "B46" "PaxDollar"
Set-TextValue "C46" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D46" "1.000"
Set-TextValue "E46" "  -0.03%  "

Set-TextValue "B47" "Aptos"
Set-TextValue "C47" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.784"
Set-TextValue "E47" "  +3.14%  "

Set-TextValue "D48" "987.83"
Set-TextValue "E48" "  +0.82%  "

Set-TextValue "D49" "0.1248"
Set-TextValue "E49" "  -1.49%  "

Set-TextValue "D50" "9.242"
Set-TextValue "E50" "  -2.95%  "

Set-TextValue "D51" "35.89"
Set-TextValue "E51" "  -0.45%  "

